function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style()
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextValue $ws "D2" "28.982.95"
Set-TextValue $ws "E2" "  -0.98%  "

# Row 3
Set-TextValue $ws "D3" "1.825.61"
Set-TextValue $ws "E3" "  -1.17%  "

# Row 4
Set-TextValue $ws "E4" "  +0.09%  "

# Row 5
Set-TextValue $ws "D5" "241.20"
Set-TextValue $ws "E5" "  +0.08%  "

# Row 6
Set-TextValue $ws "D6" "0.6362"
Set-TextValue $ws "E6" "  -5.55%  "

# Row 7
Set-TextValue $ws "D7" "1.002"
Set-TextValue $ws "E7" "  +0.21%  "

# Row 8
Set-TextValue $ws "D8" "44.41"
Set-TextValue $ws "E8" "  +5.41%  "

# Row 9
Set-TextValue $ws "D9" "0.2917"

# Row 10
Set-TextValue $ws "D10" "0.07298"
Set-TextValue $ws "E10" "  -2.15%  "

# Row 11
Set-TextValue $ws "D11" "22.70"
Set-TextValue $ws "E11" "  -1.28%  "

# Row 12
Set-TextValue $ws "D12" "0.07655"
Set-TextValue $ws "E12" "  -0.80%  "

# Row 13
Set-TextValue $ws "D13" "1.830.48"
Set-TextValue $ws "E13" "  -0.98%  "

# Row 14
Set-TextValue $ws "E14" "  -0.96%  "

# Row 15
Set-TextValue $ws "D15" "0.6611"
Set-TextValue $ws "E15" "  -1.90%  "

# Row 16
Set-TextValue $ws "D16" "81.97"
Set-TextValue $ws "E16" "  -4.86%  "

# Row 17
Set-TextValue $ws "D17" "6.040"
Set-TextValue $ws "E17" "  -2.10%  "

# Row 18
Set-TextValue $ws "D18" "0.000008571"
Set-TextValue $ws "E18" "  +2.78%  "

# Row 19
Set-TextValue $ws "D19" "29.002.31"
Set-TextValue $ws "E19" "  -0.90%  "

# Row 20
Set-TextValue $ws "D20" "2.083.21"
Set-TextValue $ws "E20" "  -0.45%  "

# Row 21
Set-TextValue $ws "D21" "12.35"
Set-TextValue $ws "E21" "  -1.55%  "

# Row 22
Set-TextValue $ws "D22" "222.89"
Set-TextValue $ws "E22" "  -2.74%  "

# Row 23
Set-TextValue $ws "D23" "1.001"
Set-TextValue $ws "E23" "  +0.04%  "

# Row 24
Set-TextValue $ws "D24" "7.061"
Set-TextValue $ws "E24" "  -2.01%  "

# Row 25
Set-TextValue $ws "E25" "  +0.10%  "

# Row 26
Set-TextValue $ws "D26" "157.85"
Set-TextValue $ws "E26" "  -2.06%  "

# Row 27
Set-TextValue $ws "D27" "8.429"
Set-TextValue $ws "E27" "  -3.43%  "

# Row 28
Set-TextValue $ws "D28" "0.1367"
Set-TextValue $ws "E28" "  -3.20%  "

# Row 29
Set-TextValue $ws "D29" "17.84"
Set-TextValue $ws "E29" "  -1.14%  "

# Row 30
Set-TextValue $ws "D30" "1.499"
Set-TextValue $ws "E30" "  -1.01%  "

# Row 31
Set-TextValue $ws "D31" "4.075"
Set-TextValue $ws "E31" "  -2.44%  "

# Row 32
Set-TextValue $ws "D32" "1.200"
Set-TextValue $ws "E32" "  +0.43%  "

# Row 33
Set-TextValue $ws "D33" "3.995"
Set-TextValue $ws "E33" "  -2.11%  "

# Row 34
Set-TextValue $ws "D34" "0.05276"
Set-TextValue $ws "E34" "  -1.14%  "

# Row 35
Set-TextValue $ws "D35" "0.7377"
Set-TextValue $ws "E35" "  -2.87%  "

# Row 36
Set-TextValue $ws "D36" "1.822"
Set-TextValue $ws "E36" "  -2.84%  "

# Row 37
Set-TextValue $ws "D37" "1.150"
Set-TextValue $ws "E37" "  +0.90%  "

# Row 38
Set-TextValue $ws "D38" "2.645"
Set-TextValue $ws "E38" "  -1.26%  "

# Row 39
Set-TextValue $ws "D39" "1.285.66"
Set-TextValue $ws "E39" "  -2.80%  "

# Row 40
Set-TextValue $ws "D40" "2.746"
Set-TextValue $ws "E40" "  +0.65%  "

# Row 41
Set-TextValue $ws "D41" "0.01780"
Set-TextValue $ws "E41" "  -1.34%  "

# Row 42
Set-TextValue $ws "D42" "6.361"
Set-TextValue $ws "E42" "  +6.03%  "

# Row 43
Set-TextValue $ws "D43" "0.8925"
Set-TextValue $ws "E43" "  -3.38%  "

# Row 44
Set-TextValue $ws "D44" "1.001"
Set-TextValue $ws "E44" "  -0.14%  "

# Row 45
Set-TextValue $ws "D45" "102.43"
Set-TextValue $ws "E45" "  -1.08%  "

# Row 46
Set-TextValue $ws "D46" "1.981.02"
Set-TextValue $ws "E46" "  -0.46%  "

# Row 47
Set-TextValue $ws "E47" "  -0.54%  "

# Row 48
Set-TextValue $ws "D48" "0.5139"
Set-TextValue $ws "E48" "  -0.55%  "

# Row 49
Set-TextValue $ws "D49" "63.85"
Set-TextValue $ws "E49" "  -0.55%  "

# Row 50
Set-TextValue $ws "B50" "RenderToken"
Set-TextValue $ws "C50" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D50" "1.722"
Set-TextValue $ws "E50" "  -3.20%  "

# Row 51
Set-TextValue $ws "B51" "XinFinNetwork"
Set-TextValue $ws "C51" "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
Set-TextValue $ws "D51" "0.07365"
Set-TextValue $ws "E51" "  -11.08%  "
